$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used row in column A (Beteckning) to bound the update.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-09 -> 2023-09-10) for every data row.
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45179
